$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-9 from 45184 to 45185
$ws.Range("C2:C9").Value = 45185
